$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "last updated" timestamp string in A1 (10:22 -> 10:52)
$ws.Range("A1").Value = "Datos actualizados a 1 de Mayo de 2020 a las 10:52"

# 2) Row 31 - Austria: refreshed case numbers
$ws.Cells.Item(31,2).Value = 15561
$ws.Cells.Item(31,3).Value = 109
$ws.Cells.Item(31,4).Value = 13110
$ws.Cells.Item(31,5).Value = 1862
$ws.Cells.Item(31,6).Value = 124
$ws.Cells.Item(31,7).Value = 5
$ws.Cells.Item(31,8).Value = 589

# 3) Row 43 - Filipinas: refreshed case numbers
$ws.Cells.Item(43,2).Value = 8772
$ws.Cells.Item(43,3).Value = 284
$ws.Cells.Item(43,4).Value = 1084
$ws.Cells.Item(43,5).Value = 7109
$ws.Cells.Item(43,6).Value = 31
$ws.Cells.Item(43,7).Value = 11
$ws.Cells.Item(43,8).Value = 579

# 4) Row 44 - Banglades is newly inserted into the ranking here, pushing
#    Noruega and Chequia down one place each (rows 45 and 46)
$ws.Cells.Item(44,1).Value = "Banglades"
$ws.Cells.Item(44,2).Value = 8238
$ws.Cells.Item(44,3).Value = 571
$ws.Cells.Item(44,4).Value = 174
$ws.Cells.Item(44,5).Value = 7894
$ws.Cells.Item(44,6).Value = 1
$ws.Cells.Item(44,7).Value = 2
$ws.Cells.Item(44,8).Value = 170

# 5) Row 45 - Noruega (formerly row 44, values unchanged, just shifted down)
$ws.Cells.Item(45,1).Value = "Noruega"
$ws.Cells.Item(45,2).Value = 7738
$ws.Cells.Item(45,3).Value = 0
$ws.Cells.Item(45,4).Value = 32
$ws.Cells.Item(45,5).Value = 7496
$ws.Cells.Item(45,6).Value = 37
$ws.Cells.Item(45,7).Value = 0
$ws.Cells.Item(45,8).Value = 210

# 6) Row 46 - Chequia (formerly row 45, values unchanged, just shifted down)
$ws.Cells.Item(46,1).Value = "Chequia"
$ws.Cells.Item(46,2).Value = 7689
$ws.Cells.Item(46,3).Value = 7
$ws.Cells.Item(46,4).Value = 3314
$ws.Cells.Item(46,5).Value = 4138
$ws.Cells.Item(46,6).Value = 67
$ws.Cells.Item(46,7).Value = 1
$ws.Cells.Item(46,8).Value = 237

# 7) Row 90 - Hong Kong: refreshed case numbers
$ws.Cells.Item(90,2).Value = 1040
$ws.Cells.Item(90,3).Value = 2
$ws.Cells.Item(90,4).Value = 859
$ws.Cells.Item(90,5).Value = 177

# 8) Row 167 - Nepal now ranks above Polinesia Francesa, with refreshed numbers
$ws.Cells.Item(167,1).Value = "Nepal"
$ws.Cells.Item(167,2).Value = 59
$ws.Cells.Item(167,3).Value = 2
$ws.Cells.Item(167,4).Value = 16
$ws.Cells.Item(167,5).Value = 43
$ws.Cells.Item(167,6).Value = 0

# 9) Row 168 - Polinesia Francesa (formerly row 167, values unchanged, shifted down)
$ws.Cells.Item(168,1).Value = "Polinesia Francesa"
$ws.Cells.Item(168,2).Value = 58
$ws.Cells.Item(168,4).Value = 50
$ws.Cells.Item(168,5).Value = 8
$ws.Cells.Item(168,6).Value = 1
